$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete "Resolving-Mac" target-cluster rows (old rows 14-17),
# which collapses the dimension from A1:T17 down to A1:T13.
$ws.Range("A14:T17").Clear()

# Rewrite the remaining data rows (2-13) with the refreshed TPM-derived values.
# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Dhh"
$ws.Cells.Item(2, 3).Value = "Boc"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 4.063713
$ws.Cells.Item(2, 8).Value = 12.191139
$ws.Cells.Item(2, 9).Value = 0.5065008440615062
$ws.Cells.Item(2, 10).Value = 0.5065008440615063
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 2.537011
$ws.Cells.Item(2, 14).Value = 7.611032999999999
$ws.Cells.Item(2, 15).Value = 0.05020703468023843
$ws.Cells.Item(2, 16).Value = 0.05020703468023844
$ws.Cells.Item(2, 17).Value = 10.309684581843
$ws.Cells.Item(2, 18).Value = 92.78716123658698
$ws.Cells.Item(2, 19).Value = 0.02542990544336608
$ws.Cells.Item(2, 20).Value = 0.02542990544336609

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Dhh"
$ws.Cells.Item(3, 3).Value = "Boc"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 4.063713
$ws.Cells.Item(3, 8).Value = 12.191139
$ws.Cells.Item(3, 9).Value = 0.5065008440615062
$ws.Cells.Item(3, 10).Value = 0.5065008440615063
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 40.76140833333334
$ws.Cells.Item(3, 14).Value = 122.284225
$ws.Cells.Item(3, 15).Value = 0.8066616352105005
$ws.Cells.Item(3, 16).Value = 0.8066616352105006
$ws.Cells.Item(3, 17).Value = 165.642664942475
$ws.Cells.Item(3, 18).Value = 1490.783984482275
$ws.Cells.Item(3, 19).Value = 0.4085747991061533
$ws.Cells.Item(3, 20).Value = 0.4085747991061535

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Dhh"
$ws.Cells.Item(4, 3).Value = "Boc"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 4.063713
$ws.Cells.Item(4, 8).Value = 12.191139
$ws.Cells.Item(4, 9).Value = 0.5065008440615062
$ws.Cells.Item(4, 10).Value = 0.5065008440615063
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 7.232567333333333
$ws.Cells.Item(4, 14).Value = 21.697702
$ws.Cells.Item(4, 15).Value = 0.143131330109261
$ws.Cells.Item(4, 16).Value = 0.143131330109261
$ws.Cells.Item(4, 17).Value = 29.391077895842
$ws.Cells.Item(4, 18).Value = 264.519701062578
$ws.Cells.Item(4, 19).Value = 0.07249613951198676
$ws.Cells.Item(4, 20).Value = 0.07249613951198679

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Dhh"
$ws.Cells.Item(5, 3).Value = "Boc"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 2.089228666666667
$ws.Cells.Item(5, 8).Value = 6.267686
$ws.Cells.Item(5, 9).Value = 0.2604012840237886
$ws.Cells.Item(5, 10).Value = 0.2604012840237886
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 2.537011
$ws.Cells.Item(5, 14).Value = 7.611032999999999
$ws.Cells.Item(5, 15).Value = 0.05020703468023843
$ws.Cells.Item(5, 16).Value = 0.05020703468023844
$ws.Cells.Item(5, 17).Value = 5.300396108848666
$ws.Cells.Item(5, 18).Value = 47.703564979638
$ws.Cells.Item(5, 19).Value = 0.01307397629776097
$ws.Cells.Item(5, 20).Value = 0.01307397629776098

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Dhh"
$ws.Cells.Item(6, 3).Value = "Boc"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 2.089228666666667
$ws.Cells.Item(6, 8).Value = 6.267686
$ws.Cells.Item(6, 9).Value = 0.2604012840237886
$ws.Cells.Item(6, 10).Value = 0.2604012840237886
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 40.76140833333334
$ws.Cells.Item(6, 14).Value = 122.284225
$ws.Cells.Item(6, 15).Value = 0.8066616352105005
$ws.Cells.Item(6, 16).Value = 0.8066616352105006
$ws.Cells.Item(6, 17).Value = 85.15990278370558
$ws.Cells.Item(6, 18).Value = 766.4391250533502
$ws.Cells.Item(6, 19).Value = 0.2100557255815433
$ws.Cells.Item(6, 20).Value = 0.2100557255815434

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Dhh"
$ws.Cells.Item(7, 3).Value = "Boc"
$ws.Cells.Item(7, 4).Value = "MuSCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 2.089228666666667
$ws.Cells.Item(7, 8).Value = 6.267686
$ws.Cells.Item(7, 9).Value = 0.2604012840237886
$ws.Cells.Item(7, 10).Value = 0.2604012840237886
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 7.232567333333333
$ws.Cells.Item(7, 14).Value = 21.697702
$ws.Cells.Item(7, 15).Value = 0.143131330109261
$ws.Cells.Item(7, 16).Value = 0.143131330109261
$ws.Cells.Item(7, 17).Value = 15.11048700639689
$ws.Cells.Item(7, 18).Value = 135.994383057572
$ws.Cells.Item(7, 19).Value = 0.03727158214448431
$ws.Cells.Item(7, 20).Value = 0.03727158214448432

# Row 8
$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 2).Value = "Dhh"
$ws.Cells.Item(8, 3).Value = "Boc"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 1.588356333333333
$ws.Cells.Item(8, 8).Value = 4.765069
$ws.Cells.Item(8, 9).Value = 0.1979725988286506
$ws.Cells.Item(8, 10).Value = 0.1979725988286507
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 2.537011
$ws.Cells.Item(8, 14).Value = 7.611032999999999
$ws.Cells.Item(8, 15).Value = 0.05020703468023843
$ws.Cells.Item(8, 16).Value = 0.05020703468023844
$ws.Cells.Item(8, 17).Value = 4.029677489586333
$ws.Cells.Item(8, 18).Value = 36.26709740627699
$ws.Cells.Item(8, 19).Value = 0.009939617135126994
$ws.Cells.Item(8, 20).Value = 0.009939617135126996

# Row 9
$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 2).Value = "Dhh"
$ws.Cells.Item(9, 3).Value = "Boc"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 1.588356333333333
$ws.Cells.Item(9, 8).Value = 4.765069
$ws.Cells.Item(9, 9).Value = 0.1979725988286506
$ws.Cells.Item(9, 10).Value = 0.1979725988286507
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 40.76140833333334
$ws.Cells.Item(9, 14).Value = 122.284225
$ws.Cells.Item(9, 15).Value = 0.8066616352105005
$ws.Cells.Item(9, 16).Value = 0.8066616352105006
$ws.Cells.Item(9, 17).Value = 64.74364108183612
$ws.Cells.Item(9, 18).Value = 582.6927697365251
$ws.Cells.Item(9, 19).Value = 0.1596969002979917
$ws.Cells.Item(9, 20).Value = 0.1596969002979918

# Row 10
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Dhh"
$ws.Cells.Item(10, 3).Value = "Boc"
$ws.Cells.Item(10, 4).Value = "MuSCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 1.588356333333333
$ws.Cells.Item(10, 8).Value = 4.765069
$ws.Cells.Item(10, 9).Value = 0.1979725988286506
$ws.Cells.Item(10, 10).Value = 0.1979725988286507
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 7.232567333333333
$ws.Cells.Item(10, 14).Value = 21.697702
$ws.Cells.Item(10, 15).Value = 0.143131330109261
$ws.Cells.Item(10, 16).Value = 0.143131330109261
$ws.Cells.Item(10, 17).Value = 11.48789413015978
$ws.Cells.Item(10, 18).Value = 103.391047171438
$ws.Cells.Item(10, 19).Value = 0.02833608139553189
$ws.Cells.Item(10, 20).Value = 0.02833608139553189

# Row 11
$ws.Cells.Item(11, 1).Value = "Resolving-Mac"
$ws.Cells.Item(11, 2).Value = "Dhh"
$ws.Cells.Item(11, 3).Value = "Boc"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.281814
$ws.Cells.Item(11, 8).Value = 0.845442
$ws.Cells.Item(11, 9).Value = 0.03512527308605438
$ws.Cells.Item(11, 10).Value = 0.03512527308605439
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 2.537011
$ws.Cells.Item(11, 14).Value = 7.611032999999999
$ws.Cells.Item(11, 15).Value = 0.05020703468023843
$ws.Cells.Item(11, 16).Value = 0.05020703468023844
$ws.Cells.Item(11, 17).Value = 0.7149652179539999
$ws.Cells.Item(11, 18).Value = 6.434686961585999
$ws.Cells.Item(11, 19).Value = 0.001763535803984378
$ws.Cells.Item(11, 20).Value = 0.001763535803984379

# Row 12
$ws.Cells.Item(12, 1).Value = "Resolving-Mac"
$ws.Cells.Item(12, 2).Value = "Dhh"
$ws.Cells.Item(12, 3).Value = "Boc"
$ws.Cells.Item(12, 4).Value = "FAPs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.281814
$ws.Cells.Item(12, 8).Value = 0.845442
$ws.Cells.Item(12, 9).Value = 0.03512527308605438
$ws.Cells.Item(12, 10).Value = 0.03512527308605439
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 40.76140833333334
$ws.Cells.Item(12, 14).Value = 122.284225
$ws.Cells.Item(12, 15).Value = 0.8066616352105005
$ws.Cells.Item(12, 16).Value = 0.8066616352105006
$ws.Cells.Item(12, 17).Value = 11.48713552805
$ws.Cells.Item(12, 18).Value = 103.38421975245
$ws.Cells.Item(12, 19).Value = 0.02833421022481201
$ws.Cells.Item(12, 20).Value = 0.02833421022481202

# Row 13
$ws.Cells.Item(13, 1).Value = "Resolving-Mac"
$ws.Cells.Item(13, 2).Value = "Dhh"
$ws.Cells.Item(13, 3).Value = "Boc"
$ws.Cells.Item(13, 4).Value = "MuSCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.281814
$ws.Cells.Item(13, 8).Value = 0.845442
$ws.Cells.Item(13, 9).Value = 0.03512527308605438
$ws.Cells.Item(13, 10).Value = 0.03512527308605439
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 7.232567333333333
$ws.Cells.Item(13, 14).Value = 21.697702
$ws.Cells.Item(13, 15).Value = 0.143131330109261
$ws.Cells.Item(13, 16).Value = 0.143131330109261
$ws.Cells.Item(13, 17).Value = 2.038238730476
$ws.Cells.Item(13, 18).Value = 18.344148574284
$ws.Cells.Item(13, 19).Value = 0.00502752705725799
$ws.Cells.Item(13, 20).Value = 0.005027527057257991

Write-Host "Updated Dhh-Boc sheet: dropped Resolving-Mac target rows and refreshed TPM values."
